$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 12120592
$ws.Range("I19").Value = 8538252
$ws.Range("J19").Value = 20001738
$ws.Range("K19").Value = 8538252
$ws.Range("L19").Value = 20001738
$ws.Range("M19").Value = -8538077
$ws.Range("N19").Value = -20002088
$ws.Range("H39").Value = 92.25
$ws.Range("I39").Value = 92.25
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 276.75
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 19.25
$ws.Range("N39").ClearContents()
$ws.Range("H64").Value = 3792.282
$ws.Range("I64").Value = 3681.818
$ws.Range("J64").Value = 3935.2354
$ws.Range("K64").Value = 3681.818
$ws.Range("L64").Value = 3935.2354
$ws.Range("M64").Value = -3433.818
$ws.Range("N64").Value = -4431.2354
$ws.Range("H67").Value = 3792.282
$ws.Range("I67").Value = 3681.818
$ws.Range("J67").Value = 3935.2354
$ws.Range("K67").Value = 3681.818
$ws.Range("L67").Value = 3935.2354
$ws.Range("M67").Value = -2823.818
$ws.Range("N67").Value = -5651.2354
$ws.Range("H82").Value = 6030080.5
$ws.Range("I82").Value = 6030080.5
$ws.Range("K82").Value = 18090241.5
$ws.Range("M82").Value = -18089835.5
$ws.Range("H85").Value = 6030080.5
$ws.Range("I85").Value = 6030080.5
$ws.Range("K85").Value = 18090241.5
$ws.Range("M85").Value = -18088837.5
$ws.Range("H112").Value = 1333.8055
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 1354.7715
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 4064.3145
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -6280.3145
$ws.Range("H132").Value = 3581.3022
$ws.Range("I132").Value = 2936.7368
$ws.Range("K132").Value = 8810.2104
$ws.Range("M132").Value = -6280.2104

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17245716
$ws.Range("I32").Value = 19610692
$ws.Range("K32").Value = 19610692
$ws.Range("M32").Value = -19610405
$ws.Range("H74").Value = 1225.0454
$ws.Range("I74").Value = 1135.5897
$ws.Range("K74").Value = 1135.5897
$ws.Range("M74").Value = -261.5897
$ws.Range("H77").Value = 1225.0454
$ws.Range("I77").Value = 1135.5897
$ws.Range("K77").Value = 5677.9485
$ws.Range("M77").Value = -1309.9485
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 2028.0588
$ws.Range("I122").Value = 1956.4166
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 5869.2498
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -3419.2498
$ws.Range("N122").Value = -11500
$ws.Range("H123").Value = 38596.332
$ws.Range("J123").Value = 38596.332
$ws.Range("L123").Value = 38596.332
$ws.Range("N123").Value = -48396.332
$ws.Range("H132").Value = 2270.5
$ws.Range("I132").Value = 1798.4
$ws.Range("K132").Value = 5395.200000000001
$ws.Range("M132").Value = -2865.200000000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2756.37
$ws.Range("I105").Value = 1382.75
$ws.Range("J105").Value = 2943.682
$ws.Range("K105").Value = 1382.75
$ws.Range("L105").Value = 2943.682
$ws.Range("M105").Value = 364.25
$ws.Range("N105").Value = -6437.682

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2200.5854
$ws.Range("I31").Value = 1650.7667
$ws.Range("J31").Value = 3700.0908
$ws.Range("K31").Value = 1650.7667
$ws.Range("L31").Value = 3700.0908
$ws.Range("M31").Value = -1355.7667
$ws.Range("N31").Value = -4290.0908
$ws.Range("H34").Value = 2200.5854
$ws.Range("I34").Value = 1650.7667
$ws.Range("J34").Value = 3700.0908
$ws.Range("K34").Value = 1650.7667
$ws.Range("L34").Value = 3700.0908
$ws.Range("M34").Value = -1448.7667
$ws.Range("N34").Value = -4104.0908
$ws.Range("H58").Value = 923.87036
$ws.Range("I58").Value = 796.87805
$ws.Range("J58").Value = 1324.3846
$ws.Range("K58").Value = 796.87805
$ws.Range("L58").Value = 1324.3846
$ws.Range("M58").Value = -593.87805
$ws.Range("N58").Value = -1730.3846
$ws.Range("H86").Value = 4615.9287
$ws.Range("I86").Value = 6086.7144
$ws.Range("J86").Value = 3145.1428
$ws.Range("K86").Value = 6086.7144
$ws.Range("L86").Value = 3145.1428
$ws.Range("M86").Value = -4963.7144
$ws.Range("N86").Value = -5391.1428
$ws.Range("H89").Value = 4615.9287
$ws.Range("I89").Value = 6086.7144
$ws.Range("J89").Value = 3145.1428
$ws.Range("K89").Value = 30433.572
$ws.Range("L89").Value = 15725.714
$ws.Range("M89").Value = -24817.572
$ws.Range("N89").Value = -26957.714
$ws.Range("H132").Value = 1443.425
$ws.Range("I132").Value = 1166.4286
$ws.Range("J132").Value = 3382.4
$ws.Range("K132").Value = 3499.2858
$ws.Range("L132").Value = 10147.2
$ws.Range("M132").Value = -969.2857999999997
$ws.Range("N132").Value = -15207.2
$ws.Range("H134").Value = 1521.619
$ws.Range("I134").Value = 894.1142599999999
$ws.Range("K134").Value = 2682.34278
$ws.Range("M134").Value = -147.3427799999999
$ws.Range("H136").Value = 923.87036
$ws.Range("I136").Value = 796.87805
$ws.Range("J136").Value = 1324.3846
$ws.Range("K136").Value = 2390.63415
$ws.Range("L136").Value = 3973.1538
$ws.Range("M136").Value = 159.3658500000001
$ws.Range("N136").Value = -9073.1538

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1916
$ws.Range("I64").Value = 12
$ws.Range("K64").Value = 36
$ws.Range("M64").Value = 234
$ws.Range("H67").Value = 1916
$ws.Range("I67").Value = 12
$ws.Range("K67").Value = 36
$ws.Range("M67").Value = 900
$ws.Range("H114").Value = 2180.7368
$ws.Range("I114").Value = 1267.3334
$ws.Range("J114").Value = 3002.8
$ws.Range("K114").Value = 3802.0002
$ws.Range("L114").Value = 9008.400000000001
$ws.Range("M114").Value = -548.0001999999999
$ws.Range("N114").Value = -15516.4
$ws.Range("H132").Value = 778461
$ws.Range("I132").Value = 1260.5
$ws.Range("K132").Value = 11344.5
$ws.Range("M132").Value = -8814.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 996.26086
$ws.Range("I102").Value = 985.7143
$ws.Range("J102").Value = 1107
$ws.Range("K102").Value = 985.7143
$ws.Range("L102").Value = 1107
$ws.Range("M102").Value = 636.2857
$ws.Range("N102").Value = -4351
$ws.Range("H113").Value = 8939.214
$ws.Range("I113").Value = 863.7143
$ws.Range("J113").Value = 17014.715
$ws.Range("K113").Value = 863.7143
$ws.Range("L113").Value = 17014.715
$ws.Range("M113").Value = 1306.2857
$ws.Range("N113").Value = -21354.715
$ws.Range("H123").Value = 14886.187
$ws.Range("J123").Value = 14886.187
$ws.Range("L123").Value = 14886.187
$ws.Range("N123").Value = -19786.187
$ws.Range("H124").Value = 31980
$ws.Range("J124").Value = 31980
$ws.Range("L124").Value = 31980
$ws.Range("N124").Value = -41800
$ws.Range("H125").Value = 27000
$ws.Range("J125").Value = 27000
$ws.Range("L125").Value = 27000
$ws.Range("N125").Value = -31920

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2185.7144
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 2433.3333
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 2433.3333
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -2657.3333
$ws.Range("H40").Value = 6254.154
$ws.Range("I40").Value = 6300.5
$ws.Range("J40").Value = 6180
$ws.Range("K40").Value = 6300.5
$ws.Range("L40").Value = 6180
$ws.Range("M40").Value = -6164.5
$ws.Range("N40").Value = -6452
$ws.Range("H122").Value = 3723.8
$ws.Range("I122").Value = 3599.5
$ws.Range("J122").Value = 3806.6667
$ws.Range("K122").Value = 10798.5
$ws.Range("L122").Value = 11420.0001
$ws.Range("M122").Value = -8348.5
$ws.Range("N122").Value = -16320.0001
$ws.Range("H126").Value = 2185.7144
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2433.3333
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 7299.999899999999
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -12239.9999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1923.6875
$ws.Range("I126").Value = 1851.9333
$ws.Range("K126").Value = 5555.7999
$ws.Range("M126").Value = -3085.7999
$ws.Range("H132").Value = 2192.3928
$ws.Range("I132").Value = 1427.8823
$ws.Range("J132").Value = 3373.9092
$ws.Range("K132").Value = 4283.6469
$ws.Range("L132").Value = 10121.7276
$ws.Range("M132").Value = -1753.6469
$ws.Range("N132").Value = -15181.7276
